# Auto-generated Excel COM-interop edit script
# Applies cell-level value corrections to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 374.8889
$ws.Range("I5").Value = 145
$ws.Range("K5").Value = 145
$ws.Range("M5").Value = -30
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 86
$ws.Range("H86").Value = 4065.7222
$ws.Range("I86").Value = 3376.1
$ws.Range("J86").Value = 4927.75
$ws.Range("K86").Value = 3376.1
$ws.Range("L86").Value = 4927.75
$ws.Range("M86").Value = -2253.1
$ws.Range("N86").Value = -7173.75
# Row 89
$ws.Range("H89").Value = 4065.7222
$ws.Range("I89").Value = 3376.1
$ws.Range("J89").Value = 4927.75
$ws.Range("K89").Value = 16880.5
$ws.Range("L89").Value = 24638.75
$ws.Range("M89").Value = -11264.5
$ws.Range("N89").Value = -35870.75
# Row 99
$ws.Range("H99").Value = 494.16666
$ws.Range("I99").Value = 453
$ws.Range("K99").Value = 1359
$ws.Range("M99").Value = 139
# Row 100
$ws.Range("H100").Value = 34443.613
$ws.Range("I100").Value = 42220.12
$ws.Range("J100").Value = 2041.5
$ws.Range("K100").Value = 42220.12
$ws.Range("L100").Value = 2041.5
$ws.Range("M100").Value = -41679.12
$ws.Range("N100").Value = -3123.5
# Row 106
$ws.Range("H106").Value = 20436.291
$ws.Range("I106").Value = 7768
$ws.Range("J106").Value = 31155.615
$ws.Range("K106").Value = 7768
$ws.Range("L106").Value = 31155.615
$ws.Range("M106").Value = -7137
$ws.Range("N106").Value = -32417.615
# Row 112
$ws.Range("H112").Value = 1740.963
$ws.Range("I112").Value = 1098
$ws.Range("K112").Value = 3294
$ws.Range("M112").Value = -2186
# Row 132
$ws.Range("H132").Value = 1861.0385
$ws.Range("I132").Value = 1795.75
$ws.Range("K132").Value = 5387.25
$ws.Range("M132").Value = -2857.25
# Row 138
$ws.Range("H138").Value = 2834.8484
$ws.Range("J138").Value = 3198.2
$ws.Range("L138").Value = 9594.599999999999
$ws.Range("N138").Value = -19874.6
# Row 141
$ws.Range("H141").Value = 2222.842
$ws.Range("J141").Value = 2044.5
$ws.Range("L141").Value = 6133.5
$ws.Range("N141").Value = -16493.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
# Row 32
$ws.Range("H32").Value = 2996.5625
$ws.Range("I32").Value = 1841.262
$ws.Range("K32").Value = 1841.262
$ws.Range("M32").Value = -1554.262
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 14399.2
$ws.Range("I122").Value = 12499.75
$ws.Range("K122").Value = 37499.25
$ws.Range("M122").Value = -35049.25
# Row 139
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
# Row 99
$ws.Range("H99").Value = 4717.375
$ws.Range("I99").Value = 5098.1
$ws.Range("J99").Value = 2813.75
$ws.Range("K99").Value = 5098.1
$ws.Range("L99").Value = 2813.75
$ws.Range("M99").Value = -3600.1
$ws.Range("N99").Value = -5809.75
# Row 134
$ws.Range("H134").Value = 6628.079
$ws.Range("I134").Value = 4509.1
$ws.Range("K134").Value = 13527.3
$ws.Range("M134").Value = -10992.3

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2643.5293
$ws.Range("I16").Value = 2527.1035
$ws.Range("K16").Value = 2527.1035
$ws.Range("M16").Value = -2240.1035
# Row 113
$ws.Range("H113").Value = 2643.5293
$ws.Range("I113").Value = 2527.1035
$ws.Range("K113").Value = 2527.1035
$ws.Range("M113").Value = -357.1035000000002
# Row 125
$ws.Range("H125").Value = 72621.57000000001
$ws.Range("J125").Value = 73110
$ws.Range("L125").Value = 73110
$ws.Range("N125").Value = -78030
# Row 132
$ws.Range("H132").Value = 24212.92
$ws.Range("I132").Value = 15062.841
$ws.Range("K132").Value = 45188.523
$ws.Range("M132").Value = -42658.523
# Row 134
$ws.Range("H134").Value = 7618.8423
$ws.Range("I134").Value = 6694.6665
$ws.Range("J134").Value = 9887.272000000001
$ws.Range("K134").Value = 20083.9995
$ws.Range("L134").Value = 29661.816
$ws.Range("M134").Value = -17548.9995
$ws.Range("N134").Value = -34731.81600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 2568.2856
$ws.Range("I86").Value = 990
$ws.Range("J86").Value = 3199.6
$ws.Range("K86").Value = 2970
$ws.Range("L86").Value = 9598.799999999999
$ws.Range("N86").Value = -11970.8
$ws.Range("M86").Value = -1784
# Row 89
$ws.Range("H89").Value = 2568.2856
$ws.Range("I89").Value = 990
$ws.Range("J89").Value = 3199.6
$ws.Range("K89").Value = 8910
$ws.Range("L89").Value = 28796.4
$ws.Range("N89").Value = -40652.39999999999
$ws.Range("M89").Value = -2982

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9973.857
$ws.Range("I70").Value = 7283
$ws.Range("K70").Value = 7283
$ws.Range("M70").Value = -7013
# Row 73
$ws.Range("H73").Value = 9973.857
$ws.Range("I73").Value = 7283
$ws.Range("K73").Value = 7283
$ws.Range("M73").Value = -6347
# Row 97
$ws.Range("H97").Value = 590.45
$ws.Range("I97").Value = 492.7647
$ws.Range("K97").Value = 492.7647
$ws.Range("M97").Value = 3.235299999999995
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
# Row 22
$ws.Range("H22").Value = 3157.1875
$ws.Range("I22").Value = 3253.182
$ws.Range("J22").Value = 2946
$ws.Range("K22").Value = 3253.182
$ws.Range("L22").Value = 2946
$ws.Range("M22").Value = -2958.182
$ws.Range("N22").Value = -3536
# Row 27
$ws.Range("H27").Value = 3157.1875
$ws.Range("I27").Value = 3253.182
$ws.Range("J27").Value = 2946
$ws.Range("K27").Value = 3253.182
$ws.Range("L27").Value = 2946
$ws.Range("M27").Value = -3146.182
$ws.Range("N27").Value = -3160
# Row 82
$ws.Range("H82").Value = 1785.0385
$ws.Range("I82").Value = 1602
$ws.Range("J82").Value = 2077.9
$ws.Range("K82").Value = 1602
$ws.Range("L82").Value = 2077.9
$ws.Range("M82").Value = -1241
$ws.Range("N82").Value = -2799.9
# Row 85
$ws.Range("H85").Value = 1785.0385
$ws.Range("I85").Value = 1602
$ws.Range("J85").Value = 2077.9
$ws.Range("K85").Value = 1602
$ws.Range("L85").Value = 2077.9
$ws.Range("M85").Value = -354
$ws.Range("N85").Value = -4573.9
# Row 93
$ws.Range("H93").Value = 5229.75
$ws.Range("I93").Value = 5723.1665
$ws.Range("K93").Value = 5723.1665
$ws.Range("M93").Value = -4475.1665
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 5558226
$ws.Range("I132").Value = 5558226
$ws.Range("K132").Value = 16674678
$ws.Range("M132").Value = -16672148
# Row 136
$ws.Range("H136").Value = 3705655.5
$ws.Range("I136").Value = 4446239.5
$ws.Range("K136").Value = 13338718.5
$ws.Range("M136").Value = -13336168.5

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 62
$ws.Range("H62").Value = 1866.1111
$ws.Range("I62").Value = 1748.25
$ws.Range("J62").Value = 1960.4
$ws.Range("K62").Value = 1748.25
$ws.Range("L62").Value = 1960.4
$ws.Range("M62").Value = -1124.25
$ws.Range("N62").Value = -3208.4
# Row 65
$ws.Range("H65").Value = 1866.1111
$ws.Range("I65").Value = 1748.25
$ws.Range("J65").Value = 1960.4
$ws.Range("K65").Value = 8741.25
$ws.Range("L65").Value = 9802
$ws.Range("M65").Value = -5621.25
$ws.Range("N65").Value = -16042
# Row 107
$ws.Range("H107").Value = 1166.6976
$ws.Range("I107").Value = 970.5714
$ws.Range("J107").Value = 1532.8
$ws.Range("K107").Value = 2911.7142
$ws.Range("L107").Value = 4598.4
$ws.Range("M107").Value = -991.7142000000003
$ws.Range("N107").Value = -8438.4
# Row 126
$ws.Range("H126").Value = 4474.2915
$ws.Range("I126").Value = 4716.095
$ws.Range("K126").Value = 14148.285
$ws.Range("M126").Value = -11678.285
# Row 132
$ws.Range("H132").Value = 16573.44
$ws.Range("I132").Value = 10527.2
$ws.Range("K132").Value = 31581.6
$ws.Range("M132").Value = -29051.6
# Row 136
$ws.Range("H136").Value = 702.2857
$ws.Range("I136").Value = 702.2857
$ws.Range("K136").Value = 2106.8571
$ws.Range("M136").Value = 443.1428999999998

